# Update latest output (run 57)
# Applies the new optimisation-run values to the "Schedule" summary sheet
# and the "Detailed" per-interval sheet.

$wb = $excel.ActiveWorkbook

# --- Schedule sheet: refresh the single summary row (Cost, Unit Cost) ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 661.6262347500001
$wsSchedule.Range("F2").Value = 10.93958721478175

# --- Detailed sheet: refresh per-interval Price values, and flip a couple
#     of intervals from "forecast" to "historical" now that they have
#     actually elapsed. ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Cells.Item(19, 2).Value = 33.22362

$wsDetailed.Cells.Item(20, 2).Value = 35.88

$wsDetailed.Cells.Item(21, 2).Value = -9.5031
$wsDetailed.Cells.Item(21, 3).Value = "historical"

$wsDetailed.Cells.Item(22, 2).Value = -7.46266
$wsDetailed.Cells.Item(22, 3).Value = "historical"

$wsDetailed.Cells.Item(23, 2).Value = 22.07

$wsDetailed.Cells.Item(24, 2).Value = 0.7

$wsDetailed.Cells.Item(25, 2).Value = 22.07

$wsDetailed.Cells.Item(26, 2).Value = 36.06046

$wsDetailed.Cells.Item(27, 2).Value = 36.06045

$wsDetailed.Cells.Item(30, 2).Value = 23.97858

$wsDetailed.Cells.Item(31, 2).Value = 23.50558

$wsDetailed.Cells.Item(32, 2).Value = 28.08235

$wsDetailed.Cells.Item(33, 2).Value = 40.54

$wsDetailed.Cells.Item(34, 2).Value = 33.2694

$wsDetailed.Cells.Item(35, 2).Value = 8.34709

$wsDetailed.Cells.Item(36, 2).Value = -0.09452000000000001

$wsDetailed.Cells.Item(37, 2).Value = -2.99905

$wsDetailed.Cells.Item(38, 2).Value = -2.89726

$wsDetailed.Cells.Item(39, 2).Value = -2.90585

$wsDetailed.Cells.Item(41, 2).Value = 13.59537

$wsDetailed.Cells.Item(42, 2).Value = 29.68591

$wsDetailed.Cells.Item(43, 2).Value = 29.86477

$wsDetailed.Cells.Item(44, 2).Value = 15.46948

$wsDetailed.Cells.Item(48, 2).Value = 57.3

$wsDetailed.Cells.Item(49, 2).Value = 57.06005
